# Update "想去人数" (F) counts and mark row 33's "最低票价" (G) as sold out
# on both the "展览" and "全部类型" sheets (they carry duplicate data).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 106
    7  = 2088
    10 = 4561
    14 = 15
    15 = 141
    19 = 3453
    21 = 552
    24 = 88
    25 = 97
    29 = 209
    30 = 15
    31 = 690
    32 = 2095
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }

    # Row 33: ticket sold out -> replace numeric price with text
    $ws.Range("G33").Value = "已售罄"
}
